$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Time"
$ws.Range("B1").Value = "AAPL"
$ws.Range("C1").Value = "AMZN"
$ws.Range("D1").Value = "WEN"
$ws.Range("E1").Value = "TSLA"
$ws.Range("F1").Value = "MSFT"

$ws.Range("B2").Value = "'123.75"
$ws.Range("C2").Value = "'3,158.00"
$ws.Range("D2").Value = "'22.02"
$ws.Range("E2").Value = "'641.76"
$ws.Range("F2").Value = "'214.29"
